{"js": "// Replace each three-digit-by-one-digit multiplication expression with its\n// new value. Every occurrence is unique in the document, so an exact-text\n// search + replace for each pair is safe and precise.\nconst replacements = [\n  [\"307\u00d74=1228\", \"884\u00d78=7072\"],\n  [\"402\u00d77=2814\", \"313\u00d75=1565\"],\n  [\"344\u00d72=688\", \"120\u00d75=600\"],\n  [\"618\u00d77=4326\", \"556\u00d76=3336\"],\n  [\"239\u00d76=1434\", \"709\u00d79=6381\"],\n  [\"185\u00d74=740\", \"288\u00d76=1728\"],\n  [\"732\u00d78=5856\", \"154\u00d76=924\"],\n  [\"647\u00d78=5176\", \"358\u00d74=1432\"],\n  [\"196\u00d74=784\", \"397\u00d77=2779\"],\n  [\"667\u00d73=2001\", \"121\u00d73=363\"],\n  [\"406\u00d74=1624\", \"949\u00d76=5694\"],\n  [\"229\u00d78=1832\", \"997\u00d74=3988\"],\n  [\"424\u00d73=1272\", \"459\u00d77=3213\"],\n  [\"541\u00d79=4869\", \"493\u00d76=2958\"],\n  [\"511\u00d73=1533\", \"192\u00d78=1536\"],\n  [\"305\u00d76=1830\", \"552\u00d76=3312\"],\n  [\"230\u00d78=1840\", \"101\u00d72=202\"],\n  [\"964\u00d72=1928\", \"848\u00d76=5088\"],\n  [\"508\u00d75=2540\", \"857\u00d72=1714\"],\n  [\"624\u00d79=5616\", \"287\u00d73=861\"],\n  [\"693\u00d78=5544\", \"582\u00d72=1164\"],\n  [\"844\u00d74=3376\", \"238\u00d79=2142\"],\n  [\"304\u00d79=2736\", \"402\u00d73=1206\"],\n  [\"895\u00d78=7160\", \"874\u00d74=3496\"],\n  [\"597\u00d73=1791\", \"261\u00d74=1044\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit-by-one-digit multiplication expression with its\n# new value. Every occurrence is unique in the document, so an exact-text\n# Find/Replace for each pair is safe and precise.\n$pairs = @(\n  @{old=\"307\u00d74=1228\"; new=\"884\u00d78=7072\"},\n  @{old=\"402\u00d77=2814\"; new=\"313\u00d75=1565\"},\n  @{old=\"344\u00d72=688\"; new=\"120\u00d75=600\"},\n  @{old=\"618\u00d77=4326\"; new=\"556\u00d76=3336\"},\n  @{old=\"239\u00d76=1434\"; new=\"709\u00d79=6381\"},\n  @{old=\"185\u00d74=740\"; new=\"288\u00d76=1728\"},\n  @{old=\"732\u00d78=5856\"; new=\"154\u00d76=924\"},\n  @{old=\"647\u00d78=5176\"; new=\"358\u00d74=1432\"},\n  @{old=\"196\u00d74=784\"; new=\"397\u00d77=2779\"},\n  @{old=\"667\u00d73=2001\"; new=\"121\u00d73=363\"},\n  @{old=\"406\u00d74=1624\"; new=\"949\u00d76=5694\"},\n  @{old=\"229\u00d78=1832\"; new=\"997\u00d74=3988\"},\n  @{old=\"424\u00d73=1272\"; new=\"459\u00d77=3213\"},\n  @{old=\"541\u00d79=4869\"; new=\"493\u00d76=2958\"},\n  @{old=\"511\u00d73=1533\"; new=\"192\u00d78=1536\"},\n  @{old=\"305\u00d76=1830\"; new=\"552\u00d76=3312\"},\n  @{old=\"230\u00d78=1840\"; new=\"101\u00d72=202\"},\n  @{old=\"964\u00d72=1928\"; new=\"848\u00d76=5088\"},\n  @{old=\"508\u00d75=2540\"; new=\"857\u00d72=1714\"},\n  @{old=\"624\u00d79=5616\"; new=\"287\u00d73=861\"},\n  @{old=\"693\u00d78=5544\"; new=\"582\u00d72=1164\"},\n  @{old=\"844\u00d74=3376\"; new=\"238\u00d79=2142\"},\n  @{old=\"304\u00d79=2736\"; new=\"402\u00d73=1206\"},\n  @{old=\"895\u00d78=7160\"; new=\"874\u00d74=3496\"},\n  @{old=\"597\u00d73=1791\"; new=\"261\u00d74=1044\"}\n)\n\n$d = $word.ActiveDocument\n\nforeach ($p in $pairs) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $found = $find.Execute($p.old, $false, $false, $false, $false, $false, $true, 1, $false, $p.new, 2)\n  if (-not $found) {\n    Write-Output \"WARNING: text not found: $($p.old)\"\n  }\n}\n"}
